# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.870.42"
$ws.Range("E2").Value = "  -0.45%  "

$ws.Range("D3").Value = "1.887.30"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7499"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3120"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.46%  "

$ws.Range("E9").Value = "  -1.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07129"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08538"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7600"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.57%  "

$ws.Range("D13").Value = "1.894.47"
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.362"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.144"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.18%  "

$ws.Range("D17").Value = "29.905.47"
$ws.Range("E17").Value = "  -0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007798"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "

$ws.Range("D21").Value = "2.147.40"
$ws.Range("E21").Value = "  -3.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.992"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1596"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.362"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "

$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.029"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.517"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.533"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.479"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.099"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05396"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.236"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7418"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.004"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.712"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01942"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("E40").Value = "  -1.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4454"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.81%  "

$ws.Range("D42").Value = "1.101.91"
$ws.Range("E42").Value = "  -4.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.075"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8580"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.665"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.864"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.071"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.34%  "

$ws.Range("D51").Value = "2.050.76"
$ws.Range("E51").Value = "  +1.85%  "
